# Add a new "Leavefor" column (F) to the holiday sample sheet with
# "Night Shift" / "Day Shift" values for the two existing holiday rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Leavefor"
$ws.Range("F2").Value = "Night Shift"
$ws.Range("F3").Value = "Day Shift"

$ws.Range("F4").Select()
